$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Title Card file type" filenames (column G) for the rows that
# got new/updated title-card art.
$ws.Range("G7").Value  = "donut_remix_megamix.png"
$ws.Range("G16").Value = "fork_lifter_2p_fever.png"
$ws.Range("G18").Value = "pirate_crew_fever.png"
$ws.Range("G20").Value = "rhythm_fighter_fever.png"

# Column H was an empty spacer column - remove it entirely so the
# "Remix/Sequel?" / "Required Games" data shifts one column to the left.
$ws.Columns.Item(8).Delete()

# Widen column G (Title Card file type) now that it holds the longer
# megamix/fever file names.
$ws.Columns.Item(7).ColumnWidth = 24
